$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 226; this shifts the existing rows
# 226..287 down to 227..288 (and grows the used range to A1:R288).
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with a new price-survey record.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R carry the same constant values used
# throughout this "Betarraga" block; D/J/K/L/M/P hold the new record's
# date/volume/price figures.
$ws.Range("A226").Value = 5
$ws.Range("B226").Value = "Macroferia Regional de Talca"
$ws.Range("C226").Value = "Maule"
$ws.Range("D226").Value = 44642
$ws.Range("E226").Value = 7
$ws.Range("F226").Value = 100114014
$ws.Range("G226").Value = "Betarraga"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 3000
$ws.Range("K226").Value = 700
$ws.Range("L226").Value = 700
$ws.Range("M226").Value = 700
$ws.Range("N226").Value = "`$/paquete 5 unidades"
$ws.Range("O226").Value = "Región del Maule"
$ws.Range("P226").Value = 140
$ws.Range("Q226").Value = 5
$ws.Range("R226").Value = "Hortaliza"
